$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so purely numeric-looking strings
# (e.g. "1.002") are not auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.974.68'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.643.98'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '215.09'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').Value = '0.5214'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.2603'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.06356'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').Value = '20.73'
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('D11').Value = '0.07672'
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('D12').Value = '1.636.79'
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = '4.419'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').Value = '1.868.32'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').Value = '0.5523'
$ws.Range('E15').Value = '  +1.80%  '
$ws.Range('D16').Value = '0.0₅8280'
$ws.Range('E16').Value = '  +3.51%  '
$ws.Range('D17').Value = '64.64'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').Value = '25.984.84'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '4.704'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '187.97'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').Value = '10.15'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '6.256'
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = '144.35'
$ws.Range('E25').Value = '  -3.73%  '
$ws.Range('D26').Value = '0.1220'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').Value = '7.391'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = '15.83'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').Value = '1.389'
$ws.Range('E29').Value = '  +1.55%  '
$ws.Range('D30').Value = '0.05930'
$ws.Range('E30').Value = '  -5.39%  '
$ws.Range('D31').Value = '1.263'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = '3.390'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('D33').Value = '3.399'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('D34').Value = '1.650'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').Value = '0.9942'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').Value = '2.397'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '2.756'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '0.5624'
$ws.Range('E38').Value = '  -6.15%  '
$ws.Range('D39').Value = '0.01604'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').Value = '5.861'
$ws.Range('E40').Value = '  -3.20%  '
$ws.Range('D41').Value = '0.8528'
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').Value = '1.028.81'
$ws.Range('E43').Value = '  -7.32%  '
$ws.Range('D44').Value = '98.93'
$ws.Range('E44').Value = '  -1.64%  '
$ws.Range('D45').Value = '1.793.87'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').Value = '55.58'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').Value = '8.045'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').Value = '0.05145'
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('D51').Value = '0.4221'
$ws.Range('E51').Value = '  -0.33%  '
